$d = $word.ActiveDocument

# Locate the paragraph that begins the "I am interested..." section and the
# paragraph that ends with "...to your team." (the block being replaced).
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*I am interested in the*") {
        $startPara = $d.Paragraphs.Item($i)
    }
    if ($ptext -like "*a valuable asset*to your team.*") {
        $endPara = $d.Paragraphs.Item($i)
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate target paragraphs for replacement"
}

$replaceRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">I’m interested in the Mechanical Engineering position recently listed by your company. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>In competition</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> for pre-capstone,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> my team’s product achieved unparalleled success </w:t>
      </w:r>
      <w:r>
        <w:t>with multiple flawless attempts at autonomous navigation through an obstacle course. My personal contributions included machining and assembly of steering and powertrain components, design and fabrication of electrical power and control systems, and all testing. These systems included a single ultrasonic ranging sensor, a servo motor, a stepper driven linear actuator, a simple syringe-hydraulic actuator, and required both 8V and 16V</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> power availability.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Taking on the role of team captain for capstone, I applied an agile and modular approach to the design and development of a complex electrohydraulic system designed to operate autonomously in extreme conditions and with very tight geometric restrictions. Two custom hydraulic actuators, and </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>all of</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> the infrastructure associated, including the autonomous digital control system, was part of this ambitious project. My technical role included design</w:t>
      </w:r>
      <w:r>
        <w:t>, sourcing, and</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> fabrication of all infrastructural components for the hydraulics, electrical power, and digital controls</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> systems</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:t>Key components included all fittings and routing for hydraulic power, brushless DC driven gear pump, stepper driven test rigs, associated power supplies, digital communications cabling etc. To control each module and test rig, I developed a control box capable of controlling each independently with plug and play compatibility.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>Much of the design principals I applied to these projects came from working directly with many such systems on an industrial scale as an intern with IBA (Ion Beam Applications). Though not part of the design process</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> here, I’ve been working in the guts of a </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>large scale</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> particle accelerator system requiring over a Mega-Watt of power at times. The entire system is at high vacuum </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>at all times</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and uses vane pumps, turbo pumps, and oil diffusion pumps. As well, we must robotically position patients to an accuracy of &lt;1mm and deliver a high-energy proton beam to 4 treatment rooms.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>I feel that my broad experience gives me a unique ability to approach</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> design from many different perspectives, as well as facilitate effective collaboration between interdisciplinary contributors. The projects above as well as some other work I have done is documented on LinkedIn and SlideShare. </w:t>
      </w:r>
    </w:p>

'@

$replaceRange.InsertXML($newXml)

Write-Host "Replacement complete"
